# Generate Report for Handoff
# Swap the GUID-named source file (and its derived handoff-package names)
# for a freshly generated one, and bump the associated timestamps.

$wb = $excel.ActiveWorkbook

$oldGuid = "660679d4-156c-41f0-96b3-25bac457ce59"
$newGuid = "855bebd0-e14f-4b76-afd8-bfd13c3e8764"

$oldHash = "14c7f8f6f0d21cd1d19dbfacae4a49f896083d9c"
$newHash = "f96f0d11bed2aa4c5b25d9175ddaedaab9b4cc71"

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = "2016-08-26 04:55:24"
foreach ($hl in $wsOverview.Hyperlinks) {
    $hl.TextToDisplay = "e2e\$newGuid.md"
}

# ---- zh-cn sheet ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("A2").Value = "$newGuid.md"
$wsZhCn.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-26 04:55:19"
foreach ($hl in $wsZhCn.Hyperlinks) {
    $hl.TextToDisplay = "$newGuid.md"
}

# ---- de-de sheet ----
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("A2").Value = "$newGuid.md"
$wsDeDe.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-26 04:55:24"
foreach ($hl in $wsDeDe.Hyperlinks) {
    $hl.TextToDisplay = "$newGuid.md"
}
